$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E3) / Correspond Handback DateTime (H3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 08:09:35"
$wsZhCn.Range("H3").Value = "2016-03-25 08:10:21"

# de-de sheet: Correspond Handoff Datetime (E3) / Correspond Handback DateTime (H3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 08:09:44"
$wsDeDe.Range("H3").Value = "2016-03-25 08:10:37"
